$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 686 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A686:N686"))
$ws.Range("A686").Value = 45193.39322135417
$ws.Range("B686").Value = 'pks5176275@naver.com'
$ws.Range("C686").Value = '소프트웨어학부'
$ws.Range("D686").Value = 20235161
$ws.Range("E686").Value = '박광성'
$ws.Range("G686").Value = 0.2
$ws.Range("H686").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I686").Value = '952만 명'
$ws.Range("J686").Value = 0.002
$ws.Range("K686").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("M686").ClearContents()
$ws.Range("N686").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 687 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A687:N687"))
$ws.Range("A687").Value = 45193.41605400463
$ws.Range("B687").Value = 'youngsoo051400@naver.com'
$ws.Range("C687").Value = '경영대학'
$ws.Range("D687").Value = 20233052
$ws.Range("E687").Value = '최영수'
$ws.Range("G687").Value = 0.2
$ws.Range("H687").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I687").Value = '952만 명'
$ws.Range("J687").Value = 0.059
$ws.Range("K687").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M687").ClearContents()
$ws.Range("N687").Value = '모름/무응답'

# Row 688 (template row 684)
$ws.Range("A684:N684").Copy($ws.Range("A688:N688"))
$ws.Range("A688").Value = 45193.42356760417
$ws.Range("B688").Value = 'jangho5636@gmail.com'
$ws.Range("C688").Value = '러시아학과'
$ws.Range("D688").Value = 20161723
$ws.Range("E688").Value = '이장호'
$ws.Range("G688").Value = 0.2
$ws.Range("H688").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I688").Value = '779만 명'
$ws.Range("J688").Value = 0.151
$ws.Range("K688").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("M688").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("N688").ClearContents()

# Row 689 (template row 648)
$ws.Range("A648:N648").Copy($ws.Range("A689:N689"))
$ws.Range("A689").Value = 45193.42733503472
$ws.Range("B689").Value = 'naturally160@naver.com'
$ws.Range("C689").Value = '식품영양학과'
$ws.Range("D689").Value = 20233852
$ws.Range("E689").Value = '홍지원'
$ws.Range("G689").Value = 0.2
$ws.Range("H689").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I689").Value = '952만 명'
$ws.Range("J689").Value = 0.374
$ws.Range("K689").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M689").ClearContents()
$ws.Range("N689").Value = '모름/무응답'

# Row 690 (template row 664)
$ws.Range("A664:N664").Copy($ws.Range("A690:N690"))
$ws.Range("A690").Value = 45193.43617369213
$ws.Range("B690").Value = 'dungunfight9@gmail.com'
$ws.Range("C690").Value = '금융재무학과'
$ws.Range("D690").Value = 20222986
$ws.Range("E690").Value = '유홍현'
$ws.Range("G690").Value = 0.2
$ws.Range("H690").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I690").Value = '952만 명'
$ws.Range("J690").Value = 0.151
$ws.Range("K690").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("M690").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("N690").ClearContents()

# Row 691 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A691:N691"))
$ws.Range("A691").Value = 45193.4380800463
$ws.Range("B691").Value = 'ryuthgud@naver.com'
$ws.Range("C691").Value = '경영대학'
$ws.Range("D691").Value = 20232947
$ws.Range("E691").Value = '류소형'
$ws.Range("G691").Value = 0.2
$ws.Range("H691").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I691").Value = '952만 명'
$ws.Range("J691").Value = 0.059
$ws.Range("K691").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M691").ClearContents()
$ws.Range("N691").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 692 (template row 685)
$ws.Range("A685:N685").Copy($ws.Range("A692:N692"))
$ws.Range("A692").Value = 45193.4482565162
$ws.Range("B692").Value = 'rdc9118@naver.com'
$ws.Range("C692").Value = '바이오메디컬'
$ws.Range("D692").Value = 20173608
$ws.Range("E692").Value = '김예찬'
$ws.Range("G692").Value = 0.15
$ws.Range("H692").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I692").Value = '166만 명'
$ws.Range("J692").Value = 0.151
$ws.Range("K692").Value = '중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다'
$ws.Range("M692").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("N692").ClearContents()

# Row 693 (template row 672)
$ws.Range("A672:N672").Copy($ws.Range("A693:N693"))
$ws.Range("A693").Value = 45193.458744421296
$ws.Range("B693").Value = 'mhkimghrhkd@naver.com'
$ws.Range("C693").Value = '체육학과'
$ws.Range("D693").Value = 20204110
$ws.Range("E693").Value = '김민혁'
$ws.Range("G693").Value = 0.2
$ws.Range("H693").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I693").Value = '166만 명'
$ws.Range("J693").Value = 0.151
$ws.Range("K693").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("M693").ClearContents()
$ws.Range("N693").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 694 (template row 685)
$ws.Range("A685:N685").Copy($ws.Range("A694:N694"))
$ws.Range("A694").Value = 45193.483619583334
$ws.Range("B694").Value = 'dndbql123@naver.com'
$ws.Range("C694").Value = '미디어스쿨'
$ws.Range("D694").Value = 20232574
$ws.Range("E694").Value = '전유비'
$ws.Range("G694").Value = 0.2
$ws.Range("H694").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I694").Value = '952만 명'
$ws.Range("J694").Value = 0.059
$ws.Range("K694").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M694").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("N694").ClearContents()

# Row 695 (template row 685)
$ws.Range("A685:N685").Copy($ws.Range("A695:N695"))
$ws.Range("A695").Value = 45193.486902476856
$ws.Range("B695").Value = 'chs1886@naver.com'
$ws.Range("C695").Value = '경영학과'
$ws.Range("D695").Value = 20193009
$ws.Range("E695").Value = '최현승'
$ws.Range("G695").Value = 0.2
$ws.Range("H695").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I695").Value = '952만 명'
$ws.Range("J695").Value = 0.059
$ws.Range("K695").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M695").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("N695").ClearContents()

# Row 696 (template row 684)
$ws.Range("A684:N684").Copy($ws.Range("A696:N696"))
$ws.Range("A696").Value = 45193.48774652778
$ws.Range("B696").Value = 'chaeyun7206@gmail.com'
$ws.Range("C696").Value = '언어청각학부'
$ws.Range("D696").Value = 20233932
$ws.Range("E696").Value = '배채윤'
$ws.Range("G696").Value = 0.15
$ws.Range("H696").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I696").Value = '779만 명'
$ws.Range("J696").Value = 0.151
$ws.Range("K696").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("M696").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("N696").ClearContents()

# Row 697 (template row 672)
$ws.Range("A672:N672").Copy($ws.Range("A697:N697"))
$ws.Range("A697").Value = 45193.493658865744
$ws.Range("B697").Value = '123qazwsx12@naver.com'
$ws.Range("C697").Value = '경영학과'
$ws.Range("D697").Value = 20233029
$ws.Range("E697").Value = '전예진'
$ws.Range("G697").Value = 0.2
$ws.Range("H697").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I697").Value = '779만 명'
$ws.Range("J697").Value = 0.151
$ws.Range("K697").Value = '중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다'
$ws.Range("M697").ClearContents()
$ws.Range("N697").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

# Row 698 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A698:N698"))
$ws.Range("A698").Value = 45193.50467974537
$ws.Range("B698").Value = 'kimyebin0628@naver.com'
$ws.Range("C698").Value = '법학과'
$ws.Range("D698").Value = 20202711
$ws.Range("E698").Value = '김예빈'
$ws.Range("G698").Value = 0.2
$ws.Range("H698").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I698").Value = '952만 명'
$ws.Range("J698").Value = 0.059
$ws.Range("K698").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M698").ClearContents()
$ws.Range("N698").Value = '모름/무응답'

# Row 699 (template row 671)
$ws.Range("A671:N671").Copy($ws.Range("A699:N699"))
$ws.Range("A699").Value = 45193.50510252315
$ws.Range("B699").Value = 'kim061806@naver.com'
$ws.Range("C699").Value = '간호학과'
$ws.Range("D699").Value = 20236227
$ws.Range("E699").Value = '김윤서'
$ws.Range("G699").Value = 0.1
$ws.Range("H699").Value = '조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다'
$ws.Range("I699").Value = '952만 명'
$ws.Range("J699").Value = 0.002
$ws.Range("K699").Value = '중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다'
$ws.Range("M699").ClearContents()
$ws.Range("N699").Value = '모름/무응답'

# Row 700 (template row 672)
$ws.Range("A672:N672").Copy($ws.Range("A700:N700"))
$ws.Range("A700").Value = 45193.511774479164
$ws.Range("B700").Value = 'yeonjoon9900@gmail.com'
$ws.Range("C700").Value = '간호학과'
$ws.Range("D700").Value = 20236275
$ws.Range("E700").Value = '이연준'
$ws.Range("G700").Value = 0.2
$ws.Range("H700").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I700").Value = '779만 명'
$ws.Range("J700").Value = 0.374
$ws.Range("K700").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("M700").ClearContents()
$ws.Range("N700").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'

# Row 701 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A701:N701"))
$ws.Range("A701").Value = 45193.5139671875
$ws.Range("B701").Value = 'jonahkim4415@gmail.com'
$ws.Range("C701").Value = '화학과'
$ws.Range("D701").Value = 20223409
$ws.Range("E701").Value = '김요나'
$ws.Range("G701").Value = 0.2
$ws.Range("H701").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I701").Value = '779만 명'
$ws.Range("J701").Value = 0.059
$ws.Range("K701").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("M701").ClearContents()
$ws.Range("N701").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 702 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A702:N702"))
$ws.Range("A702").Value = 45193.51402070602
$ws.Range("B702").Value = 'tnqls06082@naver.com'
$ws.Range("C702").Value = '금융재무학과'
$ws.Range("D702").Value = 20221728
$ws.Range("E702").Value = '조수빈'
$ws.Range("G702").Value = 0.1
$ws.Range("H702").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I702").Value = '952만 명'
$ws.Range("J702").Value = 0.059
$ws.Range("K702").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M702").ClearContents()
$ws.Range("N702").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

# Row 703 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A703:N703"))
$ws.Range("A703").Value = 45193.518680532405
$ws.Range("B703").Value = 'tjdals041122@gmail.com'
$ws.Range("C703").Value = '사회복지학부'
$ws.Range("D703").Value = 20232357
$ws.Range("E703").Value = '전성민'
$ws.Range("G703").Value = 0.2
$ws.Range("H703").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I703").Value = '952만 명'
$ws.Range("J703").Value = 0.059
$ws.Range("K703").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M703").ClearContents()
$ws.Range("N703").Value = '모름/무응답'

# Row 704 (template row 674)
$ws.Range("A674:N674").Copy($ws.Range("A704:N704"))
$ws.Range("A704").Value = 45193.53551234954
$ws.Range("B704").Value = 'ssw12099@naver.com'
$ws.Range("C704").Value = '소프트웨어'
$ws.Range("D704").Value = 20235190
$ws.Range("E704").Value = '소선웅'
$ws.Range("G704").Value = 0.15
$ws.Range("H704").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I704").Value = '38만 명'
$ws.Range("J704").Value = 0.059
$ws.Range("K704").Value = '중소기업이 신고법인수의 91%를 차지하는 데 부담하는 세액은 24.6%이다'
$ws.Range("M704").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("N704").ClearContents()

# Row 705 (template row 682)
$ws.Range("A682:N682").Copy($ws.Range("A705:N705"))
$ws.Range("A705").Value = 45193.53993293981
$ws.Range("B705").Value = 'wnsaus_0522@naver.com'
$ws.Range("C705").Value = '인문학부'
$ws.Range("D705").Value = 20231093
$ws.Range("E705").Value = '진희원'
$ws.Range("G705").Value = 0.2
$ws.Range("H705").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I705").Value = '166만 명'
$ws.Range("J705").Value = 0.002
$ws.Range("K705").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("M705").ClearContents()
$ws.Range("N705").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'

Write-Host "done"